$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New contacts appended below the existing three (rows 5-16, column A = Name).
$names = @(
    "David Brown",
    "Olivia Martinez",
    "Christopher Taylor",
    "Sophia Anderson",
    "Daniel Thompson",
    "Ava Garcia",
    "Matthew Rodriguez",
    "Isabella Clark",
    "James Wilson",
    "Charlotte Thomas",
    "Alexander White",
    "Mia Davis"
)

# Matching phone numbers for rows 5-9 (column B). Rows 10-16 all end up
# with the same number as row 10 (source-data bug being replicated here).
$phones = @(
    "0834567890",
    "0645678901",
    "0606789012",
    "0767890123",
    "0748901234",
    "0819012345"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

for ($i = 0; $i -lt $phones.Count; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 2).Value = $phones[$i]
}

# Rows 10-16 (column B) all share the last phone number from the list above.
for ($row = 11; $row -le 16; $row++) {
    $ws.Cells.Item($row, 2).Value = "0819012345"
}

# The saved selection moves from H14 to B1.
$ws.Range("B1").Select() | Out-Null
